$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A250").Value = 'Iñupiat Heritage Center'
$ws.Range("A251").Value = 'Ice Age Floods National Geologic Trail'
$ws.Range("A252").Value = 'Ice Age National Scenic Trail'
$ws.Range("A253").Value = 'Independence National Historical Park'
$ws.Range("A254").Value = 'Indiana Dunes National Park'
$ws.Range("A255").Value = 'Isle Royale National Park'
$ws.Range("A257").Value = 'James A Garfield National Historic Site'
$ws.Range("A258").Value = 'Jean Lafitte National Historical Park and Preserve'
$ws.Range("A259").Value = 'Jewel Cave National Monument'
$ws.Range("A260").Value = 'Jimmy Carter National Historic Site'
$ws.Range("A261").Value = 'John Day Fossil Beds National Monument'
$ws.Range("A262").Value = 'John Fitzgerald Kennedy National Historic Site'
$ws.Range("A263").Value = 'John H. Chafee Blackstone River Valley National Heritage Corridor'
$ws.Range("A264").Value = 'John Muir National Historic Site'
$ws.Range("A265").Value = 'Johnstown Flood National Memorial'
$ws.Range("A266").Value = 'Joshua Tree National Park'
$ws.Range("A267").Value = 'Journey Through Hallowed Ground National Heritage Area'
$ws.Range("A268").Value = 'Juan Bautista de Anza National Historic Trail'
$ws.Range("A270").Value = 'Kalaupapa National Historical Park'
$ws.Range("A271").Value = 'Kaloko-Honokōhau National Historical Park'
$ws.Range("A272").Value = 'Katahdin Woods and Waters National Monument'
$ws.Range("A273").Value = 'Katmai National Park & Preserve'
$ws.Range("A274").Value = 'Kenai Fjords National Park'
$ws.Range("A275").Value = 'Kenilworth Park & Aquatic Gardens'
$ws.Range("A276").Value = 'Kennesaw Mountain National Battlefield Park'
$ws.Range("A277").Value = 'Keweenaw National Historical Park'
$ws.Range("A278").Value = 'Kings Mountain National Military Park'
$ws.Range("A279").Value = 'Klondike Gold Rush - Seattle Unit National Historical Park'
$ws.Range("A280").Value = 'Klondike Gold Rush National Historical Park'
$ws.Range("A281").Value = 'Knife River Indian Villages National Historic Site'
$ws.Range("A282").Value = 'Kobuk Valley National Park'
$ws.Range("A283").Value = 'Korean War Veterans Memorial'
$ws.Range("A285").Value = 'Lake Clark National Park & Preserve'
$ws.Range("A286").Value = 'Lake Mead National Recreation Area'
$ws.Range("A287").Value = 'Lake Meredith National Recreation Area'
$ws.Range("A288").Value = 'Lake Roosevelt National Recreation Area'
$ws.Range("A289").Value = 'Lassen Volcanic National Park'
$ws.Range("A290").Value = 'Lava Beds National Monument'
$ws.Range("A291").Value = 'LBJ Memorial Grove on the Potomac'
$ws.Range("A292").Value = 'Lewis & Clark National Historic Trail'
$ws.Range("A293").Value = 'Lewis and Clark National Historical Park'
$ws.Range("A294").Value = 'Lincoln Boyhood National Memorial'
$ws.Range("A295").Value = 'Lincoln Home National Historic Site'
$ws.Range("A296").Value = 'Lincoln Memorial'
$ws.Range("A297").Value = 'Little Bighorn Battlefield National Monument'
$ws.Range("A298").Value = 'Little River Canyon National Preserve'
$ws.Range("A299").Value = 'Little Rock Central High School National Historic Site'
$ws.Range("A300").Value = 'Longfellow House Washington''s Headquarters National Historic Site'
$ws.Range("A301").Value = 'Lowell National Historical Park'
$ws.Range("A302").Value = 'Lower Delaware National Wild and Scenic River'
$ws.Range("A303").Value = 'Lower East Side Tenement Museum National Historic Site'
$ws.Range("A304").Value = 'Lyndon B Johnson National Historical Park'
$ws.Range("A306").Value = 'Maggie L Walker National Historic Site'
$ws.Range("A307").Value = 'Maine Acadian Culture'
$ws.Range("A308").Value = 'Mammoth Cave National Park'
$ws.Range("A309").Value = 'Manassas National Battlefield Park'
$ws.Range("A310").Value = 'Manhattan Project National Historical Park'
$ws.Range("A311").Value = 'Manzanar National Historic Site'
$ws.Range("A312").Value = 'Marsh - Billings - Rockefeller National Historical Park'
$ws.Range("A313").Value = 'Martin Luther King, Jr. Memorial'
$ws.Range("A314").Value = 'Martin Luther King, Jr. National Historical Park'
$ws.Range("A315").Value = 'Martin Van Buren National Historic Site'
$ws.Range("A316").Value = 'Mary McLeod Bethune Council House National Historic Site'
$ws.Range("A317").Value = 'Mesa Verde National Park'
$ws.Range("A318").Value = 'Minidoka National Historic Site'
$ws.Range("A319").Value = 'Minute Man National Historical Park'
$ws.Range("A320").Value = 'Minuteman Missile National Historic Site'
$ws.Range("A321").Value = 'Mississippi Delta National Heritage Area'
$ws.Range("A322").Value = 'Mississippi Gulf National Heritage Area'
$ws.Range("A323").Value = 'Mississippi Hills National Heritage Area'
$ws.Range("A324").Value = 'Mississippi National River and Recreation Area'
$ws.Range("A325").Value = 'Missouri National Recreational River'
$ws.Range("A326").Value = 'Mojave National Preserve'
$ws.Range("A327").Value = 'Monocacy National Battlefield'
$ws.Range("A328").Value = 'Montezuma Castle National Monument'
$ws.Range("A329").Value = 'Moores Creek National Battlefield'
$ws.Range("A330").Value = 'Mormon Pioneer National Historic Trail'
$ws.Range("A331").Value = 'Morristown National Historical Park'
$ws.Range("A332").Value = 'Motor Cities National Heritage Area'
$ws.Range("A333").Value = 'Mount Rainier National Park'
$ws.Range("A334").Value = 'Mount Rushmore National Memorial'
$ws.Range("A335").Value = 'Muir Woods National Monument'
$ws.Range("A336").Value = 'Muscle Shoals National Heritage Area'
$ws.Range("A338").Value = 'Natchez National Historical Park'
$ws.Range("A339").Value = 'Natchez Trace National Scenic Trail'
$ws.Range("A340").Value = 'Natchez Trace Parkway'
$ws.Range("A341").Value = 'National Aviation Heritage Area'
$ws.Range("A342").Value = 'National Capital Parks-East'
$ws.Range("A343").Value = 'National Mall and Memorial Parks'
$ws.Range("A344").Value = 'National Park of American Samoa'
$ws.Range("A345").Value = 'National Parks of New York Harbor'
$ws.Range("A346").Value = 'Natural Bridges National Monument'
$ws.Range("A347").Value = 'Navajo National Monument'
$ws.Range("A348").Value = 'New Bedford Whaling National Historical Park'
$ws.Range("A349").Value = 'New England National Scenic Trail'
$ws.Range("A350").Value = 'New Jersey Pinelands National Reserve'
$ws.Range("A351").Value = 'New Orleans Jazz National Historical Park'
$ws.Range("A352").Value = 'New River Gorge National River'
$ws.Range("A353").Value = 'Nez Perce National Historical Park'
$ws.Range("A354").Value = 'Niagara Falls National Heritage Area'
$ws.Range("A355").Value = 'Nicodemus National Historic Site'
$ws.Range("A356").Value = 'Ninety Six National Historic Site'
$ws.Range("A357").Value = 'Niobrara National Scenic River'
$ws.Range("A358").Value = 'Noatak National Preserve'
$ws.Range("A359").Value = 'North Cascades National Park'
$ws.Range("A360").Value = 'North Country National Scenic Trail'

$ws.Range("A360").Select()
